$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of marks for Janifer
$ws.Range("A8").Value = "Janifer"
$ws.Range("B8").Value = 23
$ws.Range("C8").Value = 45
$ws.Range("D8").Value = 53
$ws.Range("E8").Value = 121
$ws.Range("F8").Value = "C"

# Match the formatting of the other "C" grade rows (e.g. F4, F6)
$ws.Range("F4").Copy()
$ws.Range("F8").PasteSpecial(-4122)

# Update selection to D9 (the cell below the new row, matching the diff)
$ws.Range("D9").Select()
